# PowerUnit.xlsx update
# - Row 5 (C9..C40 group): Manufacturer/PartNum/Digikey ref now sourced from Samsung
# - New "Quantity*2" column (G) = Quantity * 2 for every part row
# - Selection moved to K12 (cosmetic, matches author's last cursor position)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the C9..C40 designator row's manufacturer info ---
$ws.Range("C5").Value = "Samsung"
$ws.Range("D5").Value = "CL10B683KB8NNNC"
$ws.Range("E5").Value = "1276-1814-1-ND "

# Re-apply the original text-cell formatting (setting .Value resets the
# style to a generic numeric style) by pasting formats from a cell in the
# same row that kept its original style.
$ws.Range("A5").Copy()
$ws.Range("C5:E5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 2. Add the "Quantity*2" column ---
$ws.Range("G1").Value = "Quantity*2"
for ($r = 2; $r -le 38; $r++) {
    $qty = $ws.Cells.Item($r, 6).Value2
    $ws.Cells.Item($r, 7).Value = $qty * 2
}

# Copy header/data formatting from column F onto the new column G so the
# new cells share the same border/font/number styles as the rest of the
# table.
$ws.Range("F1:F38").Copy()
$ws.Range("G1:G38").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match column width of the other columns (character width 18).
$ws.Columns.Item(7).ColumnWidth = 17.17

# --- 3. Cosmetic: move the active selection ---
$ws.Range("K12").Select()
